$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A28").Value() = 'JPA'
$ws.Range("B28").Value() = 'JPA Basics:'
$ws.Range("C28").Value() = '1. We use annoation @Entity(name="EntityName") to define a class as entity which its instance is a row in the relational database tables.
2. @Table(name=”TABLE_NAME”) > JPA will create a database table with name TABLE_NAME, but @Entity also creates a database table, so @Table is optional.
3. @Transient any field within an entity will not be persisted so no database table column will be created
4. @Version for concurrent modification and optimistic locking
5. @Id we can define the primary key
6. @GeneratedValue we can specify that the database is going to generate the value for the given field. Usually set strategyType = AUTO for id.
7. The EntityManager is like the role of contextManager in Spring. We manipulate entities with the manager, we have method: persist(), remove(), merge(), delete() from the manager.
8. JPA support SQL, but it also has its own JPQL, which is Entity-Oriented.'

$ws.Range("A29").Value() = 'JPA'
$ws.Range("B29").Value() = 'JPA vs. Hibernate'
$ws.Range("C29").Value() = '1. JPA is the standard specification and Hibernate/EclipseLink are the vendor who actually produce the implementing tools. 
2. In JPA, we config the persistence.xml to specify which vendor we are going to use. 
3. Sometimes we use Hibernate directly because JPA does not support specific feature yet. How? Change some import path from points at JPA to points at Hibernate.'

$ws.Range("A30").Value() = 'Spring-boot'
$ws.Range("B30").Value() = 'Spring Boot Concept'
$ws.Range("C30").Value() = '1. Spring is design to fasten the building of Spring application and its of course base on spring (btw, spring is build upon the vision that allow people to build software at a easier way with higher quality)
2. Spring boot also package tomcat, jboss, jetty into itself thus enabling an awesome simple way to build micro-service.'

$ws.Range("A31").Value() = 'Spring-boot'
$ws.Range("B31").Value() = 'Spring web example (case from udemy course)'
$ws.Range("C31").Value() = '* pom:
  - specifies the dependency of spring-boot-starter-web and spring-boot-starter-tomcat
  - the packaging tag with war value. 
  - build>plugin>spring-boot-maven-plugin>executable>true
* Student Class:
  - defined with @Component annoation
  - has a private Address member to which we put @Autowired annoation
* Address Class:
  - defined with @Component annoation
* In App class(the entrance class)
  - @EnableAutoConfiguration, @RestController, @ComponentScan, 3 annoations are attached to Class App.
  - its private member Student has an @Autowired
  - it has hello method with @RequestMapping("/index")
* $mvn package && mvn spring-boot:run 
'

$ws.Range("A32").Value() = 'JSR'
$ws.Range("B32").Value() = 'JSR Concepts'
$ws.Range("C32").Value() = '1. JSR stands for Java Specification Request, it defines the "Floorplan" of the java language.
2. JCP (Java Community Process) has 4 major steps that renew the JSR:
 i) Initiation: A specification is initiated by community members and approved for development by the Executive Committee. At times, there are new JSRs being accepted every week. 
 ii) Draft Releases: Once a JSR is approved, a group of experts is formed to develop a progressive drafts of the specification that anyone with an internet connection can review. 
 iii) Final Release: The Expert Group uses the public feedback to further revise the document into a Proposed Final Draft...Once approved, the final Specification, Reference Implementation and Technology Compatibility Kit are published, and the Specification Lead arranges for a Maintenance Lead.
 iv) Maintenance: The Maintenance Lead tracks requests for clarification, interpretation, enhancements and revisions in an Issue Tracker... until the specification can be revised by an Expert Group in a new JSR. 
3. Currently the EC(Executive Committee) include big names like: Azul, Eclipse Foundation, Fujitsu, Goldman Sachs, HP, IBM, Intel, Oracle, RedHat, SAP, Twitter...'

$ws.Range("A33").Value() = 'JSR'
$ws.Range("B33").Value() = 'JSR Concepts (2)'
$ws.Range("C33").Value() = 'Oracle為Java 提供實現JSR的lib或其他工具, 但如果其他Vendor有意, 也可制作可滿足該JSR的lib. 這就像JPA定義了接口, Hibernate來實現, 又例如JDBC的接口, 各大DB Vendor各自提供JDBC Connector'

$ws.Range("A34").Value() = 'Encryption'
$ws.Range("B34").Value() = 'Popular Libs'
$ws.Range("C34").Value() = 'JaSypt amd bouncy castle '

# New rows keep the sheet's standard 33pt row height, except row 31
# (Spring web example / pom) which is taller to fit its longer notes.
$ws.Rows.Item(28).RowHeight() = 33
$ws.Rows.Item(29).RowHeight() = 33
$ws.Rows.Item(30).RowHeight() = 33
$ws.Rows.Item(31).RowHeight() = 71.4
$ws.Rows.Item(32).RowHeight() = 33
$ws.Rows.Item(33).RowHeight() = 33
$ws.Rows.Item(34).RowHeight() = 33

# Match the final selection/scroll position recorded after the edits
$ws.Application.ActiveWindow.ScrollRow = 26
$ws.Range("C35").Select()